# Auto-generated edit script applying the Golem_Profits market-data refresh
# (static numeric cell updates across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 875
$ws.Range("I32").Value = 750
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 750
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -424
$ws.Range("N32").Value = -1652
$ws.Range("H69").Value = 2112.6667
$ws.Range("I69").Value = 1911.5
$ws.Range("K69").Value = 5734.5
$ws.Range("M69").Value = -4860.5
$ws.Range("H72").Value = 2112.6667
$ws.Range("I72").Value = 1911.5
$ws.Range("K72").Value = 17203.5
$ws.Range("M72").Value = -12835.5
$ws.Range("H92").Value = 55556276
$ws.Range("J92").Value = 782.5714
$ws.Range("L92").Value = 782.5714
$ws.Range("N92").Value = -3278.5714
$ws.Range("H107").Value = 70131.16
$ws.Range("I107").Value = 90851.39999999999
$ws.Range("J107").Value = 1063.6666
$ws.Range("K107").Value = 90851.39999999999
$ws.Range("L107").Value = 1063.6666
$ws.Range("M107").Value = -88931.39999999999
$ws.Range("N107").Value = -4903.6666
$ws.Range("H115").Value = 185.25
$ws.Range("I115").Value = 188
$ws.Range("J115").Value = 177
$ws.Range("K115").Value = 564
$ws.Range("L115").Value = 531
$ws.Range("M115").Value = 1003
$ws.Range("N115").Value = -3665
$ws.Range("H118").Value = 200
$ws.Range("I118").Value = 200
$ws.Range("K118").Value = 600
$ws.Range("M118").Value = 1057
$ws.Range("H125").Value = 1461
$ws.Range("I125").Value = 1139.4
$ws.Range("J125").Value = 2265
$ws.Range("K125").Value = 10254.6
$ws.Range("L125").Value = 20385
$ws.Range("M125").Value = -7794.6
$ws.Range("N125").Value = -25305
$ws.Range("H131").Value = 495
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("H135").Value = 683.4286
$ws.Range("I135").Value = 683.4286
$ws.Range("K135").Value = 6150.8574
$ws.Range("M135").Value = -3615.8574
$ws.Range("H137").Value = 1633.5
$ws.Range("I137").Value = 952.75
$ws.Range("J137").Value = 2995
$ws.Range("K137").Value = 2858.25
$ws.Range("L137").Value = 8985
$ws.Range("M137").Value = -308.25
$ws.Range("N137").Value = -14085
$ws.Range("H138").Value = 6339.269
$ws.Range("J138").Value = 6646.1
$ws.Range("L138").Value = 19938.3
$ws.Range("N138").Value = -30218.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 3000
$ws.Range("J20").Value = 3000
$ws.Range("L20").Value = 3000
$ws.Range("N20").Value = -3494
$ws.Range("I45").Value = 2499.8333
$ws.Range("K45").Value = 2499.8333
$ws.Range("M45").Value = -2122.8333
$ws.Range("H96").Value = 25998.5
$ws.Range("J96").Value = 25998.5
$ws.Range("L96").Value = 25998.5
$ws.Range("N96").Value = -31490.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 12482.1
$ws.Range("I36").Value = 10727.75
$ws.Range("K36").Value = 10727.75
$ws.Range("M36").Value = -10193.75
$ws.Range("H86").Value = 1911.2778
$ws.Range("I86").Value = 1934.25
$ws.Range("K86").Value = 1934.25
$ws.Range("M86").Value = -811.25
$ws.Range("H88").Value = 21633.572
$ws.Range("J88").Value = 21633.572
$ws.Range("L88").Value = 21633.572
$ws.Range("N88").Value = -22445.572
$ws.Range("H89").Value = 1911.2778
$ws.Range("I89").Value = 1934.25
$ws.Range("K89").Value = 9671.25
$ws.Range("M89").Value = -4055.25
$ws.Range("H91").Value = 21633.572
$ws.Range("J91").Value = 21633.572
$ws.Range("L91").Value = 21633.572
$ws.Range("N91").Value = -24441.572
$ws.Range("H107").Value = 34198.383
$ws.Range("I107").Value = 46398
$ws.Range("J107").Value = 6749.25
$ws.Range("K107").Value = 46398
$ws.Range("L107").Value = 6749.25
$ws.Range("M107").Value = -44478
$ws.Range("N107").Value = -10589.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2605.8333
$ws.Range("I11").Value = 495.33334
$ws.Range("J11").Value = 4716.3335
$ws.Range("K11").Value = 495.33334
$ws.Range("L11").Value = 4716.3335
$ws.Range("M11").Value = -355.33334
$ws.Range("N11").Value = -4996.3335
$ws.Range("H13").Value = 487.5
$ws.Range("J13").Value = 487.5
$ws.Range("L13").Value = 487.5
$ws.Range("N13").Value = -765.5
$ws.Range("H17").Value = 1009
$ws.Range("J17").Value = 1009
$ws.Range("L17").Value = 1009
$ws.Range("N17").Value = -1357
$ws.Range("H141").Value = 709081
$ws.Range("J141").Value = 709081
$ws.Range("L141").Value = 709081
$ws.Range("N141").Value = -719441

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1649.1666
$ws.Range("I4").Value = 406.8
$ws.Range("K4").Value = 1220.4
$ws.Range("M4").Value = -1108.4
$ws.Range("H32").Value = 7035
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H46").Value = 2746.2307
$ws.Range("I46").Value = 1003
$ws.Range("J46").Value = 2891.5
$ws.Range("K46").Value = 3009
$ws.Range("L46").Value = 8674.5
$ws.Range("M46").Value = -2918
$ws.Range("N46").Value = -8856.5
$ws.Range("H80").Value = 2127.2
$ws.Range("I80").Value = 1926.6666
$ws.Range("J80").Value = 2428
$ws.Range("K80").Value = 5779.9998
$ws.Range("L80").Value = 7284
$ws.Range("M80").Value = -4843.9998
$ws.Range("N80").Value = -9156
$ws.Range("H83").Value = 2127.2
$ws.Range("I83").Value = 1926.6666
$ws.Range("J83").Value = 2428
$ws.Range("K83").Value = 17339.9994
$ws.Range("L83").Value = 21852
$ws.Range("M83").Value = -12659.9994
$ws.Range("N83").Value = -31212
$ws.Range("H107").Value = 609.6667
$ws.Range("I107").Value = 298.33334
$ws.Range("J107").Value = 765.3333
$ws.Range("K107").Value = 895.0000200000001
$ws.Range("L107").Value = 2295.9999
$ws.Range("M107").Value = 1024.99998
$ws.Range("N107").Value = -6135.9999
$ws.Range("H122").Value = 1075
$ws.Range("I122").Value = 900
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 8100
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -5650
$ws.Range("N122").Value = -16150
$ws.Range("H123").Value = 2000
$ws.Range("I123").Value = 2000
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 6000
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -3550
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H134").Value = 1286.6666
$ws.Range("I134").Value = 1286.6666
$ws.Range("K134").Value = 3859.9998
$ws.Range("M134").Value = 1210.0002
$ws.Range("H138").Value = 530
$ws.Range("I138").Value = 530
$ws.Range("K138").Value = 1590
$ws.Range("M138").Value = 3550
$ws.Range("H140").Value = 759.5714
$ws.Range("I140").Value = 552.8333
$ws.Range("K140").Value = 1658.4999
$ws.Range("M140").Value = 3521.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 53000
$ws.Range("J74").Value = 53000
$ws.Range("L74").Value = 53000
$ws.Range("N74").Value = -54872
$ws.Range("H77").Value = 53000
$ws.Range("J77").Value = 53000
$ws.Range("L77").Value = 159000
$ws.Range("N77").Value = -168360
$ws.Range("H97").Value = 2637.1667
$ws.Range("I97").Value = 2878
$ws.Range("J97").Value = 2300
$ws.Range("K97").Value = 2878
$ws.Range("L97").Value = 2300
$ws.Range("M97").Value = -2382
$ws.Range("N97").Value = -3292
$ws.Range("H107").Value = 37038570
$ws.Range("I107").Value = 943.75
$ws.Range("K107").Value = 943.75
$ws.Range("M107").Value = 976.25
$ws.Range("I113").Value = 3648.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3648.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1478.75
$ws.Range("N113").ClearContents()
$ws.Range("H123").Value = 75000
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 33366.5
$ws.Range("J2").Value = 39839.8
$ws.Range("L2").Value = 39839.8
$ws.Range("N2").Value = -40063.8
$ws.Range("H68").Value = 7037.5
$ws.Range("I68").Value = 4325
$ws.Range("K68").Value = 4325
$ws.Range("M68").Value = -3576
$ws.Range("H71").Value = 7037.5
$ws.Range("I71").Value = 4325
$ws.Range("K71").Value = 21625
$ws.Range("M71").Value = -17881
$ws.Range("H93").Value = 33340164
$ws.Range("I93").Value = 41672970
$ws.Range("J93").Value = 8949.5
$ws.Range("K93").Value = 41672970
$ws.Range("L93").Value = 8949.5
$ws.Range("M93").Value = -41671722
$ws.Range("N93").Value = -11445.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 30559.6
$ws.Range("I44").Value = 5000
$ws.Range("J44").Value = 36949.5
$ws.Range("K44").Value = 5000
$ws.Range("L44").Value = 36949.5
$ws.Range("M44").Value = -4446
$ws.Range("N44").Value = -38057.5
$ws.Range("H113").Value = 411.5
$ws.Range("I113").Value = 523
$ws.Range("K113").Value = 1569
$ws.Range("M113").Value = 601
